# Aggiunte posizione nemici, aggiornata tabella exp.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio2")

# Update raw "Exp"/"Exp opt" source values in the Tabella135 table.
$ws.Range("B16").Value = 200
$ws.Range("B24").Value = 150
$ws.Range("B25").Value = 200
$ws.Range("B26").Value = 250
$ws.Range("C27").Value = 600
$ws.Range("B28").Value = 600

# Move the active selection to reflect the author's last interaction.
$ws.Activate()
$ws.Range("F13:F30").Select()
